# Correction in SA algorithm and 746 logs
# Update fitness values (column C) in the log sheet to reflect the
# corrected simulated-annealing run output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C56").Value  = 7808
$ws.Range("C57:C63").Value = 7678
$ws.Range("C64:C68").Value = 7676
$ws.Range("C69:C85").Value = 7569
